$d = $word.ActiveDocument

# --- Change 1: remove the stray "_GoBack" bookmark that currently sits
#     between " in this list" and " without problems." in the first
#     paragraph. Word recreates this bookmark automatically at the last
#     edit point; deleting it by name removes both bookmarkStart/bookmarkEnd.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Change 2: append two new list paragraphs at the very end of the
#     document (after the "occupancyMod_full.stan" bullet, before the
#     section break) describing the new R scripts, and re-create the
#     "_GoBack" bookmark at the end of the new content (this is where
#     Word would naturally leave it after the last edit).
$contentEnd = $d.Content.End
$insertionPoint = $d.Range($contentEnd, $contentEnd)

$newParagraphsXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="14"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>bird_occupancy_</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>full.R</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="14"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>format_for_analysis.R</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$insertionPoint.InsertXML($newParagraphsXml)
